# Update the dSF column (F) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = -5
    3  = 2
    5  = 2
    6  = -2
    7  = -4
    8  = -1
    9  = 1
    10 = -4
    11 = 1
    12 = -2
    14 = -3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
